$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared "DANSKIN" surname (used by both B3 and B4) is renamed to "TEST" -
# this is placeholder/test data, per the commit message ("not for sale list").
$ws.Range("B3:B4").Value = "TEST"

# Grow the sheet's selection from just A7 to cover A7 and A4:A6 as well.
$ws.Range("A4:A7").Select()
